# Update vm_pu.xlsx results for "case with 380 kV done":
# Column B (slack bus, bus 0) setpoint changes from 1.05 to 1.02 p.u.,
# and the resulting voltage magnitudes for buses 1-4 and 6-11 (columns
# C,D,E,F,I,J,K,L,M) are updated for every data row (rows 2-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns in the sheet that hold per-bus voltage magnitude results
# (G is always 1 / unchanged, H has no data, A is the time step index).
$cols = @("B","C","D","E","F","I","J","K","L","M")

# Each row: @(rowNumber, B, C, D, E, F, I, J, K, L, M)
$data = @(
    @(2, 1.02,1.036329779833897,1.043767832648233,1.035267632337253,1.052771044941338,1.040289004610983,1.04143868769633,1.046540769436517,1.038064725921372,1.055518856626008),
    @(3, 1.02,1.037154938967108,1.044400786081594,1.035966107966737,1.053541225477026,1.040476963330802,1.041908392460722,1.046985290444099,1.038572856107397,1.056102058745902),
    @(4, 1.02,1.037689606571001,1.044810997294741,1.036419076330319,1.054040581820184,1.040597794981727,1.042212384957955,1.04727288220906,1.038901981688915,1.0564797638206),
    @(5, 1.02,1.037914554885054,1.044983603355492,1.036609743612426,1.054250748094393,1.040648402667038,1.042340196920782,1.047393774112007,1.039040423958178,1.056638629389678),
    @(6, 1.02,1.037952334851699,1.04501259360582,1.036641771469761,1.054286049766595,1.040656888759721,1.042361657865094,1.04741407167004,1.039063673548611,1.056665308163794),
    @(7, 1.02,1.037692611659696,1.044813303065491,1.036421623097678,1.054043389144092,1.040598471951383,1.04221409273624,1.047274497620176,1.038903831257127,1.056481886284332),
    @(8, 1.02,1.036608493078939,1.043981606924517,1.035503475596544,1.053031123168483,1.040352689299129,1.041597412948472,1.046691005483519,1.03823638140808,1.055715882284151),
    @(9, 1.02,1.034703841611125,1.042521097776186,1.033893386682369,1.051255116269124,1.039913570626624,1.04051127982868,1.045662548614263,1.037062851852264,1.054368719974086),
    @(10, 1.02,1.033438015363137,1.041550929263659,1.032825355892337,1.050076435780487,1.039616824245299,1.039787631877623,1.04497681155588,1.03628233557423,1.053472483469317),
    @(11, 1.02,1.032890855920898,1.041131690445305,1.03236418282699,1.049567343092942,1.039487390564838,1.039474405965693,1.044679873643831,1.035944817502022,1.053084868361036),
    @(12, 1.02,1.032687761332046,1.040976096195779,1.032193078604478,1.049378438400345,1.039439172597501,1.039358079088584,1.044569577693724,1.035819517459186,1.052940961866453),
    @(13, 1.02,1.032731319229591,1.041009465789698,1.032229772140377,1.049418950239847,1.039449521862316,1.039383030708818,1.044593236526433,1.035846391588396,1.052971827056796),
    @(14, 1.02,1.032874065094661,1.041118826318827,1.032350035289038,1.049551724181886,1.039483407715036,1.039464789950635,1.044670756537858,1.035934458741148,1.053072971546651),
    @(15, 1.02,1.032962034773703,1.041186224196214,1.032424159413631,1.049633556492636,1.039504267308033,1.039515167099819,1.044718519218783,1.035988729034557,1.053135299466776),
    @(16, 1.02,1.033474348738584,1.041578770848103,1.032855989805623,1.050110249814422,1.039625394581997,1.039808422262776,1.04499651827197,1.036304745158966,1.053498218071798),
    @(17, 1.02,1.03379596582223,1.041825234070545,1.033127212723403,1.050409612140233,1.039701123254762,1.03999240615956,1.045170898247718,1.036503095426426,1.053725992020902),
    @(18, 1.02,1.033983651327081,1.04196907383487,1.03328553692295,1.050584348839724,1.039745203708234,1.040099732187604,1.045272610151843,1.036618833160317,1.053858893096572),
    @(19, 1.02,1.03404766273307,1.042018133316489,1.033339542433821,1.050643950458792,1.039760218579051,1.040136329470019,1.045307291072937,1.036658304048909,1.0539042164049),
    @(20, 1.02,1.03376144985475,1.041798782410295,1.033098100183343,1.050377480583538,1.039693007671198,1.039972665246375,1.045152189018471,1.036481809822248,1.053701549428442),
    @(21, 1.02,1.03283202598376,1.041086618772514,1.032314615363111,1.049512620165227,1.039473433049264,1.039440713351359,1.044647928815068,1.035908523237491,1.053043185041146),
    @(22, 1.02,1.032248498918409,1.040639604667064,1.031823142263964,1.048969977351415,1.039334565057663,1.039106365694239,1.044330881170958,1.035548476301123,1.052629656944994),
    @(23, 1.02,1.032557757559984,1.040876503431027,1.032083573215774,1.049257534822366,1.039408258385586,1.039283598672706,1.044498953586446,1.035739305476252,1.052848836410999),
    @(24, 1.02,1.03377704585162,1.041810734524052,1.033111254502767,1.050391999068935,1.039696675034737,1.039981585280467,1.045160642917068,1.036491427736713,1.053712593846942),
    @(25, 1.02,1.035195552771244,1.042898064928201,1.034308696590069,1.051713327248885,1.040027801816837,1.040791999496988,1.045928452833287,1.037365920003747,1.054716671167548)
)

foreach ($rowEntry in $data) {
    $rowNum = $rowEntry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $value = $rowEntry[$i + 1]
        $ws.Range("$colLetter$rowNum").Value = $value
    }
}
